# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This script updates the DAMSLTag (column I) and DialogAct (column J) values
# for a set of rows in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = @("b",  "Acknowledge (Backchannel)")
    6   = @("ba", "Appreciation")
    17  = @("aa", "Agree/Accept")
    18  = @("aa", "Agree/Accept")
    28  = @("ba", "Appreciation")
    45  = @("aa", "Agree/Accept")
    86  = @("b",  "Acknowledge (Backchannel)")
    92  = @("aa", "Agree/Accept")
    93  = @("sd", "Statement-non-opinion")
    95  = @("sv", "Statement-opinion")
    105 = @("sd", "Statement-non-opinion")
    108 = @("ba", "Appreciation")
    117 = @("ba", "Appreciation")
    127 = @("sd", "Statement-non-opinion")
    144 = @("sd", "Statement-non-opinion")
    166 = @("ba", "Appreciation")
    170 = @("sd", "Statement-non-opinion")
    175 = @("sd", "Statement-non-opinion")
    176 = @("sd", "Statement-non-opinion")
    189 = @("b",  "Acknowledge (Backchannel)")
    192 = @("b",  "Acknowledge (Backchannel)")
    238 = @("aa", "Agree/Accept")
    250 = @("b",  "Acknowledge (Backchannel)")
    255 = @("aa", "Agree/Accept")
    268 = @("ba", "Appreciation")
    278 = @("sd", "Statement-non-opinion")
    293 = @("sd", "Statement-non-opinion")
    316 = @("sd", "Statement-non-opinion")
    330 = @("ba", "Appreciation")
    337 = @("%",  "Uninterpretable")
    338 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
